$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cell F1 "Modelo" and copy the header formatting (bold, border, centered)
# from the existing header cell A1 so the new column matches the rest of the header row.
$ws.Range("F1").Value = "Modelo"
$ws.Range("A1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats

# Update existing numeric metric values in row 2
$ws.Range("B2").Value = 0.03256334092595459
$ws.Range("C2").Value = 0.999664376588444
$ws.Range("D2").Value = 0.1291041796288097

# Add new model description cell F2
$ws.Range("F2").Value = "Pipeline(steps=[('model', RandomForestRegressor(max_depth=5))])"
